# Daily scrape update: refresh AIESEC Global Talent opportunity listings
# - replaces rows 2-5 with current data
# - appends 3 new opportunities as rows 6-8
# - widens/narrows a few columns to fit the refreshed content
# - highlights the new premium ("Yes") opportunity in column E with a yellow fill

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Resize columns C, D, F, G, H (A, B, E stay the same, set for completeness) ---
$ws.Columns.Item(1).ColumnWidth = 16.166666666666668
$ws.Columns.Item(2).ColumnWidth = 54.166666666666664
$ws.Columns.Item(3).ColumnWidth = 65.16666666666667
$ws.Columns.Item(4).ColumnWidth = 48.166666666666664
$ws.Columns.Item(5).ColumnWidth = 9.166666666666666
$ws.Columns.Item(6).ColumnWidth = 16.166666666666668
$ws.Columns.Item(7).ColumnWidth = 15.166666666666666
$ws.Columns.Item(8).ColumnWidth = 44.166666666666664

# --- Write the refreshed opportunity data into rows 2-8 ---
# Row 2
$ws.Range("A2").Value = "1328588"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1328588"
$ws.Range("C2").Value = "My Way​ Operations & Innovation Coordinator​ 2026-2027"
$ws.Range("D2").Value = "40 Düsseldorf, Germany"
$ws.Range("E2").Value = "Yes"
$ws.Range("F2").Value = "136 applicants"
$ws.Range("G2").Value = "6 - 18 Months"
$ws.Range("H2").Value = "PwC Global Partnership"

# Row 3
$ws.Range("A3").Value = "1331124"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1331124"
$ws.Range("C3").Value = "MARKETING ASSISTANT ANIMAL HEALTH"
$ws.Range("D3").Value = "Santiago de Chile, Región Metropolitana, Chile"
$ws.Range("E3").Value = "No"
$ws.Range("F3").Value = "7 applicants"
$ws.Range("G3").Value = "6 - 18 Months"
$ws.Range("H3").Value = "Boehringer Ingelheim in Chile"

# Row 4
$ws.Range("A4").Value = "1330185"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1330185"
$ws.Range("C4").Value = "Front-End Web Developer"
$ws.Range("D4").Value = "Ciudad Juárez, Chihuahua, Mexico"
$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = "35 applicants"
$ws.Range("G4").Value = "6 - 18 Months"
$ws.Range("H4").Value = "EP&O Corporation"

# Row 5
$ws.Range("A5").Value = "1329810"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1329810"
$ws.Range("C5").Value = "[IMPACT FORTALEZA] Business Development - Marketplace Specialst"
$ws.Range("D5").Value = "Fortaleza - Zone 1, Fortaleza - Ceará, Brasil"
$ws.Range("E5").Value = "No"
$ws.Range("F5").Value = "29 applicants"
$ws.Range("G5").Value = "6 - 18 Months"
$ws.Range("H5").Value = "MAKRO MOVEIS E EQUIPAMENTOS MODULADOS LTDA"

# Row 6
$ws.Range("A6").Value = "1325417"
$ws.Range("B6").Value = "https://aiesec.org/opportunity/global-talent/1325417"
$ws.Range("C6").Value = "Junior Software Engineer – AI & Internal Tools (EU ONLY)"
$ws.Range("D6").Value = "Brussels, Belgium"
$ws.Range("E6").Value = "No"
$ws.Range("F6").Value = "155 applicants"
$ws.Range("G6").Value = "6 - 18 Months"
$ws.Range("H6").Value = "Eureka Resource Mining"

# Row 7
$ws.Range("A7").Value = "1325033"
$ws.Range("B7").Value = "https://aiesec.org/opportunity/global-talent/1325033"
$ws.Range("C7").Value = "Junior Full-Stack Developer – AI & Web Projects (EU ONLY)"
$ws.Range("D7").Value = "Brussels, Belgium"
$ws.Range("E7").Value = "No"
$ws.Range("F7").Value = "170 applicants"
$ws.Range("G7").Value = "6 - 18 Months"
$ws.Range("H7").Value = "Eureka Resource Mining"

# Row 8
$ws.Range("A8").Value = "1307150"
$ws.Range("B8").Value = "https://aiesec.org/opportunity/global-talent/1307150"
$ws.Range("C8").Value = "ENGINEER"
$ws.Range("D8").Value = "Gaziantep, Türkiye"
$ws.Range("E8").Value = "No"
$ws.Range("F8").Value = "46 applicants"
$ws.Range("G8").Value = "6 - 18 Months"
$ws.Range("H8").Value = "Göymen Makarna"

# --- Highlight the PREMIUM cell of the new top row (E2 = "Yes") in yellow ---
$ws.Range("E2").Interior.Color = 65535
